# Auto update Excel log
# Appends newly-logged sensor readings (2026-01-28, ~16:52-16:53) to the
# PIR, Humidity and Temperature sheets, matching the source system's export.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PIR sheet: rows 184-196 (Date, Time, Hour, Location, Value, Status)
# ---------------------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")

# Pre-format the new Date column as Text so "2026-01-28" is stored literally
# instead of being auto-converted into a date serial number.
$wsPIR.Range("A184:A196").NumberFormat = "@"

$pirTimes = @(
    @(184, "16:52:23"),
    @(185, "16:52:26"),
    @(186, "16:52:28"),
    @(187, "16:52:33"),
    @(188, "16:52:39"),
    @(189, "16:52:44"),
    @(190, "16:52:49"),
    @(191, "16:52:54"),
    @(192, "16:52:59"),
    @(193, "16:53:04"),
    @(194, "16:53:09"),
    @(195, "16:53:14"),
    @(196, "16:53:19")
)

foreach ($entry in $pirTimes) {
    $r = $entry[0]
    $time = $entry[1]

    $wsPIR.Cells.Item($r, 1).Value = "2026-01-28"
    $wsPIR.Cells.Item($r, 2).Value = $time
    $wsPIR.Cells.Item($r, 3).Value = "16:00"
    $wsPIR.Cells.Item($r, 4).Value = "Bathroom"
    $wsPIR.Cells.Item($r, 5).Value = "No Motion"
    $wsPIR.Cells.Item($r, 6).Value = "Inactive"
}

# ---------------------------------------------------------------------------
# Humidity sheet: rows 178-193 (Date, Time, Hour, Location, Value, Status)
# ---------------------------------------------------------------------------
$wsHum = $wb.Worksheets.Item("Humidity")

$wsHum.Range("A178:A193").NumberFormat = "@"
# "87.8%"-style readings must stay literal text, not become numeric percentages.
$wsHum.Range("E178:E193").NumberFormat = "@"

$humRows = @(
    @(178, "16:52:24", "87.8%"),
    @(179, "16:52:25", "86.9%"),
    @(180, "16:52:27", "87.8%"),
    @(181, "16:52:31", "87.8%"),
    @(182, "16:52:35", "86.9%"),
    @(183, "16:52:39", "87.8%"),
    @(184, "16:52:43", "86.8%"),
    @(185, "16:52:47", "87.7%"),
    @(186, "16:52:51", "87.8%"),
    @(187, "16:52:55", "86.8%"),
    @(188, "16:52:59", "87.8%"),
    @(189, "16:53:03", "86.9%"),
    @(190, "16:53:07", "87.8%"),
    @(191, "16:53:11", "87.8%"),
    @(192, "16:53:15", "86.9%"),
    @(193, "16:53:19", "87.8%")
)

foreach ($entry in $humRows) {
    $r = $entry[0]
    $time = $entry[1]
    $value = $entry[2]

    $wsHum.Cells.Item($r, 1).Value = "2026-01-28"
    $wsHum.Cells.Item($r, 2).Value = $time
    $wsHum.Cells.Item($r, 3).Value = "16:00"
    $wsHum.Cells.Item($r, 4).Value = "Bathroom"
    $wsHum.Cells.Item($r, 5).Value = $value
    $wsHum.Cells.Item($r, 6).Value = "Active"
}

# ---------------------------------------------------------------------------
# Temperature sheet: rows 178-193 (Date, Time, Hour, Location, Value, Status)
# ---------------------------------------------------------------------------
$wsTemp = $wb.Worksheets.Item("Temperature")

$wsTemp.Range("A178:A193").NumberFormat = "@"

$tempTimes = @(
    @(178, "16:52:24"),
    @(179, "16:52:25"),
    @(180, "16:52:27"),
    @(181, "16:52:31"),
    @(182, "16:52:35"),
    @(183, "16:52:40"),
    @(184, "16:52:43"),
    @(185, "16:52:47"),
    @(186, "16:52:51"),
    @(187, "16:52:55"),
    @(188, "16:53:00"),
    @(189, "16:53:04"),
    @(190, "16:53:08"),
    @(191, "16:53:12"),
    @(192, "16:53:16"),
    @(193, "16:53:20")
)

foreach ($entry in $tempTimes) {
    $r = $entry[0]
    $time = $entry[1]

    $wsTemp.Cells.Item($r, 1).Value = "2026-01-28"
    $wsTemp.Cells.Item($r, 2).Value = $time
    $wsTemp.Cells.Item($r, 3).Value = "16:00"
    $wsTemp.Cells.Item($r, 4).Value = "Bathroom"
    $wsTemp.Cells.Item($r, 5).Value = "22.8C"
    $wsTemp.Cells.Item($r, 6).Value = "Active"
}
